$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.264.51"
$ws.Range("E2").Value = "  -1.27%  "
$ws.Range("D3").Value = "2.997.84"
$ws.Range("E3").Value = "  -2.20%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "501.04"
$ws.Range("E5").Value = "  -4.96%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "136.80"
$ws.Range("E6").Value = "  -3.86%  "
$ws.Range("E8").Value = "  -4.00%  "
$ws.Range("E9").Value = "  -5.14%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.107"
$ws.Range("E10").Value = "  -4.80%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.356"
$ws.Range("E11").Value = "  -3.91%  "
$ws.Range("D12").Value = "3.506.75"
$ws.Range("E12").Value = "  -2.28%  "
$ws.Range("E13").Value = "  -2.56%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "25.99"
$ws.Range("E14").Value = "  -4.72%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000159"
$ws.Range("E15").Value = "  -6.63%  "
$ws.Range("D16").Value = "57.278.59"
$ws.Range("E16").Value = "  -1.14%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.07"
$ws.Range("E17").Value = "  -2.61%  "
$ws.Range("D18").Value = "2.998.61"
$ws.Range("E18").Value = "  -2.23%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.62"
$ws.Range("E19").Value = "  -4.07%  "
$ws.Range("E20").Value = "  -4.12%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "320.26"
$ws.Range("E21").Value = "  -5.91%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.999"
$ws.Range("E22").Value = "  +0.08%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.73"
$ws.Range("E23").Value = "  +0.64%  "
$ws.Range("E24").Value = "  -2.39%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "63.22"
$ws.Range("E25").Value = "  -2.73%  "
$ws.Range("E26").Value = "  +0.12%  "
$ws.Range("E27").Value = "  -5.15%  "
$ws.Range("D28").Value = "0.0₃0888"
$ws.Range("E28").Value = "  -9.38%  "
$ws.Range("E29").Value = "  -5.19%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.06"
$ws.Range("E30").Value = "  -4.01%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.77"
$ws.Range("E31").Value = "  -4.37%  "
$ws.Range("E32").Value = "  -6.81%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "20.13"
$ws.Range("E33").Value = "  -4.76%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "155.14"
$ws.Range("E34").Value = "  -0.98%  "
$ws.Range("E35").Value = "  -5.00%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.76"
$ws.Range("E36").Value = "  -4.21%  "
$ws.Range("E37").Value = "  -7.40%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "24.27"
$ws.Range("E38").Value = "  -8.62%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0660"
$ws.Range("E39").Value = "  -6.28%  "
$ws.Range("D40").Value = "3.027.39"
$ws.Range("E40").Value = "  -2.34%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "37.85"
$ws.Range("E41").Value = "  -0.05%  "
$ws.Range("E42").Value = "  +0.06%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.72"
$ws.Range("E43").Value = "  -4.88%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.645"
$ws.Range("E44").Value = "  -3.19%  "
$ws.Range("D45").Value = "2.177.76"
$ws.Range("E45").Value = "  -6.61%  "
$ws.Range("E46").Value = "  -7.03%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "5.93"
$ws.Range("E47").Value = "  -1.87%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.934"
$ws.Range("E48").Value = "  -9.45%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0233"
$ws.Range("E49").Value = "  -5.10%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "19.20"
$ws.Range("E50").Value = "  -5.28%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.76"
$ws.Range("E51").Value = "  -12.80%  "
